$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "81.544.53"
$ws.Range("E2").Value = "  +5.46%  "

# Row 3
$ws.Range("D3").Value = "3.181.18"
$ws.Range("E3").Value = "  +1.26%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "208.23"
$ws.Range("E5").Value = "  +3.20%  "

# Row 6
$ws.Range("D6").Value = "635.82"
$ws.Range("E6").Value = "  +1.43%  "

# Row 7
$ws.Range("D7").Value = "0.295"
$ws.Range("E7").Value = "  +29.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  +3.56%  "

# Row 10
$ws.Range("D10").Value = "3.180.76"
$ws.Range("E10").Value = "  +1.28%  "

# Row 11
$ws.Range("E11").Value = "  +9.73%  "

# Row 12
$ws.Range("E12").Value = "  +18.01%  "

# Row 13
$ws.Range("E13").Value = "  +2.31%  "

# Row 14
$ws.Range("D14").Value = "5.38"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").Value = "3.770.08"
$ws.Range("E15").Value = "  +1.38%  "

# Row 16
$ws.Range("D16").Value = "31.99"
$ws.Range("E16").Value = "  +4.24%  "

# Row 17
$ws.Range("D17").Value = "81.523.17"
$ws.Range("E17").Value = "  +5.57%  "

# Row 18
$ws.Range("D18").Value = "3.189.18"
$ws.Range("E18").Value = "  +1.77%  "

# Row 19
$ws.Range("D19").Value = "3.26"
$ws.Range("E19").Value = "  +15.33%  "

# Row 20
$ws.Range("D20").Value = "14.25"
$ws.Range("E20").Value = "  +1.52%  "

# Row 21
$ws.Range("D21").Value = "9.24"
$ws.Range("E21").Value = "  -1.77%  "

# Row 22
$ws.Range("D22").Value = "440.75"
$ws.Range("E22").Value = "  +3.34%  "

# Row 23
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").Value = "  +8.90%  "

# Row 24
$ws.Range("E24").Value = "  +4.03%  "

# Row 25
$ws.Range("D25").Value = "5.14"
$ws.Range("E25").Value = "  +10.34%  "

# Row 26
$ws.Range("D26").Value = "11.33"
$ws.Range("E26").Value = "  +5.78%  "

# Row 27
$ws.Range("D27").Value = "3.349.61"
$ws.Range("E27").Value = "  +1.34%  "

# Row 28
$ws.Range("D28").Value = "77.13"
$ws.Range("E28").Value = "  +2.17%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0000126"
$ws.Range("E30").Value = "  +10.42%  "

# Row 31
$ws.Range("D31").Value = "9.18"
$ws.Range("E31").Value = "  +4.63%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.47%  "

# Row 33
$ws.Range("D33").Value = "569.05"
$ws.Range("E33").Value = "  +9.71%  "

# Row 34
$ws.Range("D34").Value = "1.51"
$ws.Range("E34").Value = "  +1.85%  "

# Row 35
$ws.Range("E35").Value = "  +4.36%  "

# Row 36
$ws.Range("D36").Value = "0.153"
$ws.Range("E36").Value = "  +13.53%  "

# Row 37
$ws.Range("D37").Value = "0.139"
$ws.Range("E37").Value = "  +28.03%  "

# Row 38
$ws.Range("D38").Value = "23.29"
$ws.Range("E38").Value = "  +4.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40
$ws.Range("D40").Value = "0.417"
$ws.Range("E40").Value = "  +5.45%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "5.98"
$ws.Range("E41").Value = "  +11.69%  "

# Row 42
$ws.Range("D42").Value = "3.07"
$ws.Range("E42").Value = "  +21.58%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "2.04"
$ws.Range("E43").Value = "  +16.57%  "

# Row 44
$ws.Range("D44").Value = "20.77"
$ws.Range("E44").Value = "  +3.52%  "

# Row 45
$ws.Range("D45").Value = "160.05"
$ws.Range("E45").Value = "  -2.13%  "

# Row 46
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "189.30"

# Row 48
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Value = "1.35"
$ws.Range("E48").Value = "  +5.55%  "

# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "44.78"
$ws.Range("E49").Value = "  +4.77%  "

# Row 50
$ws.Range("D50").Value = "0.788"
$ws.Range("E50").Value = "  -0.79%  "

# Row 51
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.30"
$ws.Range("E51").Value = "  +5.36%  "
